# aggiornamento 15, 16, 17 marzo
# Appends three new daily rows (227-229) to the bottom of the single data
# sheet, one per new date (Excel serial dates 44301, 44302, 44303), with
# per-municipality counts in columns B:AX and a running total in AP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 227 -> date 44301
$row227 = @(44301,1,2,0,6,12,1,6,4,1,1,7,4,6,2,0,0,1,0,1,5,59,0,2,2,4,1,1,0,1,7,1,20,1,2,0,3,7,8,0,3,184,0,0,0,0,0,2,0,0)
for ($i = 0; $i -lt $row227.Length; $i++) {
    $ws.Cells.Item(227, $i + 1).Value = $row227[$i]
}
# Match the date-column formatting (border, bold, centered, date number
# format) used by every prior row in column A.
$ws.Range("A226").Copy()
$ws.Range("A227").PasteSpecial($xlPasteFormats)

# Row 228 -> date 44302
$row228 = @(44302,3,1,1,26,14,2,6,2,3,0,1,9,7,0,0,0,3,2,4,10,27,5,0,2,8,0,2,6,2,9,0,25,3,8,0,7,2,10,1,0,213,1,0,0,0,0,0,1,0)
for ($i = 0; $i -lt $row228.Length; $i++) {
    $ws.Cells.Item(228, $i + 1).Value = $row228[$i]
}
$ws.Range("A227").Copy()
$ws.Range("A228").PasteSpecial($xlPasteFormats)

# Row 229 -> date 44303
$row229 = @(44303,4,0,0,14,8,7,4,0,2,0,1,7,8,0,0,0,1,0,1,12,51,3,2,3,5,0,1,2,1,6,1,25,1,7,0,3,5,11,0,6,204,1,0,1,0,0,0,0,0)
for ($i = 0; $i -lt $row229.Length; $i++) {
    $ws.Cells.Item(229, $i + 1).Value = $row229[$i]
}
$ws.Range("A228").Copy()
$ws.Range("A229").PasteSpecial($xlPasteFormats)
